$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F2 from the numeric year 1998 to the text value "1984"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1984"

# Update the active cell selection from F6 to G5
$ws.Range("G5").Select()
